$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ANNOUNCEMENT text in B2
$ws.Range("B2").Value = "We will be having song practice on 9/23/2025 from 6:30pm - 8:30pm and on 9/25/2025: 6:00PM to 8:00PM"

# Row height change (60 -> 30)
$ws.Rows(2).RowHeight = 30

# Selection change (C8 -> F4)
$ws.Range("F4").Select() | Out-Null
